$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Clear the values in C82:C161 while keeping cell formatting (style stays as-is)
$ws.Range("C82:C161").ClearContents()

# Update the selected cell/range on the sheet view
$ws.Activate()
$ws.Range("I19").Select()
